# Updated symbol list on Thu Dec 29 15:09:09 UTC 2022 with GitHub Actions
# Applies refreshed Price (D), Volume(1h) (E) and Hora (G) values for rows 2-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values, indexed by row 2..51. Empty string = no change.
$dValues = @("246.08","24.19","5.366","0.05737","","3.136","0.8189","0.8774","0.1378","0.06986","0.03126","0.02938","0.09404","3.737","0.001527","0.04710","0.0005988","0.006159","0.001240","0.003906","","3.534","2.139","","0.1313","","0.0002331","","","","","","","","","","","","0.03720","","","0.002798","0.007597","0.00005267","0.00000000749","0.3899","0.002758","","","")

# New "Volume(1h)" (column E) values, indexed by row 2..51. Empty string = no change.
$eValues = @("","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","","42CEJICEJIBestin24h","43LocalTradersLCT","","","","","","","")

$startRow = 2
for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $startRow + $i

    $dVal = $dValues[$i]
    if ($dVal -ne "") {
        $cell = $ws.Cells.Item($row, 4)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.Style = $origStyle
    }

    $eVal = $eValues[$i]
    if ($eVal -ne "") {
        $cell = $ws.Cells.Item($row, 5)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $eVal
        $cell.Style = $origStyle
    }

    $cell = $ws.Cells.Item($row, 7)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = "15"
    $cell.Style = $origStyle
}
